$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns that must stay text even though some look numeric (account numbers in column F)
# Use a leading apostrophe via .Formula to force text / quote-prefix, preserving leading zeros.

# Row 2: Sri Tulsi Trust- Tulsi Books
$ws.Cells.Item(2, 1).Value = 4240035910
$ws.Cells.Item(2, 2).Value = 11800
$ws.Cells.Item(2, 3).Value = 'O'
$ws.Cells.Item(2, 4).Value = 'A'
$ws.Cells.Item(2, 5).Value = 'Sri Tulsi Trust- Tulsi Books'
$ws.Cells.Item(2, 6).Formula = "'" + '910010014548454'
$ws.Cells.Item(2, 8).Value = 'UTIB0001004'
$ws.Cells.Item(2, 9).Value = 'Chowpatty'
$ws.Cells.Item(2, 10).Value = 'Chowpatty'
$ws.Cells.Item(2, 11).Value = 'India'
$ws.Cells.Item(2, 12).Value = 'O|A|Sri Tulsi Trust- Tulsi Books|910010014548454||UTIB0001004|Chowpatty|Chowpatty|India'

# Row 3: TATA AIG GENERAL INSURANCE CO
$ws.Cells.Item(3, 1).Value = 4240035907
$ws.Cells.Item(3, 2).Value = 10209
$ws.Cells.Item(3, 3).Value = 'O'
$ws.Cells.Item(3, 4).Value = 'A'
$ws.Cells.Item(3, 5).Value = 'TATA AIG GENERAL INSURANCE CO'
$ws.Cells.Item(3, 6).Formula = "'" + '0005922018'
$ws.Cells.Item(3, 8).Value = 'DEUT0784BBY'
$ws.Cells.Item(3, 9).Value = 'Mumbai'
$ws.Cells.Item(3, 10).Value = 'Mumbai'
$ws.Cells.Item(3, 11).Value = 'India'
$ws.Cells.Item(3, 12).Value = 'O|A|TATA AIG GENERAL INSURANCE CO|0005922018||DEUT0784BBY|Mumbai|Mumbai|India'

# Row 4: Mukund Dhanawade
$ws.Cells.Item(4, 1).Value = 4240035886
$ws.Cells.Item(4, 2).Value = 1805
$ws.Cells.Item(4, 3).Value = 'O'
$ws.Cells.Item(4, 4).Value = 'A'
$ws.Cells.Item(4, 5).Value = 'Mukund Dhanawade'
$ws.Cells.Item(4, 6).Formula = "'" + '695010085299'
$ws.Cells.Item(4, 8).Value = 'KKBK0000638'
$ws.Cells.Item(4, 9).Value = 'Mumbai Lower Parel'
$ws.Cells.Item(4, 10).Value = 'Mumbai Lower Parel'
$ws.Cells.Item(4, 11).Value = 'India'
$ws.Cells.Item(4, 12).Value = 'O|A|Mukund Dhanawade|695010085299||KKBK0000638|Mumbai Lower Parel|Mumbai Lower Parel|India'

# Row 5: Shekhar Laxman Chavan
$ws.Cells.Item(5, 1).Value = 4240035815
$ws.Cells.Item(5, 2).Value = 2247
$ws.Cells.Item(5, 3).Value = 'O'
$ws.Cells.Item(5, 4).Value = 'A'
$ws.Cells.Item(5, 5).Value = 'Shekhar Laxman Chavan'
$ws.Cells.Item(5, 6).Formula = "'" + '41678100000806'
$ws.Cells.Item(5, 8).Value = 'BARB0WATHAR'
$ws.Cells.Item(5, 9).Value = 'Wathar '
$ws.Cells.Item(5, 10).Value = 'Wathar '
$ws.Cells.Item(5, 11).Value = 'India'
$ws.Cells.Item(5, 12).Value = 'O|A|Shekhar Laxman Chavan|41678100000806||BARB0WATHAR|Wathar |Wathar |India'

# Row 6: Vaibhav Chavan
$ws.Cells.Item(6, 1).Value = 4240035814
$ws.Cells.Item(6, 2).Value = 955
$ws.Cells.Item(6, 3).Value = 'O'
$ws.Cells.Item(6, 4).Value = 'A'
$ws.Cells.Item(6, 5).Value = 'Vaibhav Chavan'
$ws.Cells.Item(6, 6).Formula = "'" + '41678100000494'
$ws.Cells.Item(6, 8).Value = 'BARB0WATHAR'
$ws.Cells.Item(6, 9).Value = 'Wathar '
$ws.Cells.Item(6, 10).Value = 'Wathar '
$ws.Cells.Item(6, 11).Value = 'India'
$ws.Cells.Item(6, 12).Value = 'O|A|Vaibhav Chavan|41678100000494||BARB0WATHAR|Wathar |Wathar |India'

# Row 7: Ranjana  Ravindra Bhoir
$ws.Cells.Item(7, 1).Value = 4240035810
$ws.Cells.Item(7, 2).Value = 16500
$ws.Cells.Item(7, 3).Value = 'O'
$ws.Cells.Item(7, 4).Value = 'A'
$ws.Cells.Item(7, 5).Value = 'Ranjana  Ravindra Bhoir'
$ws.Cells.Item(7, 6).Formula = "'" + '55605007149'
$ws.Cells.Item(7, 8).Value = 'MAHG0005605'
$ws.Cells.Item(7, 9).Value = 'Alonde'
$ws.Cells.Item(7, 10).Value = 'Alonde'
$ws.Cells.Item(7, 11).Value = 'India'
$ws.Cells.Item(7, 12).Value = 'O|A|Ranjana  Ravindra Bhoir|55605007149||MAHG0005605|Alonde|Alonde|India'

# Row 8: Hemant Kondu Patil
$ws.Cells.Item(8, 1).Value = 4240035803
$ws.Cells.Item(8, 2).Value = 248
$ws.Cells.Item(8, 3).Value = 'O'
$ws.Cells.Item(8, 4).Value = 'A'
$ws.Cells.Item(8, 5).Value = 'Hemant Kondu Patil'
$ws.Cells.Item(8, 6).Formula = "'" + '009210100011231'
$ws.Cells.Item(8, 8).Value = 'BKID0000092'
$ws.Cells.Item(8, 9).Value = 'Boiser'
$ws.Cells.Item(8, 10).Value = 'Boiser'
$ws.Cells.Item(8, 11).Value = 'India'
$ws.Cells.Item(8, 12).Value = 'O|A|Hemant Kondu Patil|009210100011231||BKID0000092|Boiser|Boiser|India'

# Row 9: Laxman Sudhakar Padwale
$ws.Cells.Item(9, 1).Value = 4240035799
$ws.Cells.Item(9, 2).Value = 6435
$ws.Cells.Item(9, 3).Value = 'O'
$ws.Cells.Item(9, 4).Value = 'A'
$ws.Cells.Item(9, 5).Value = 'Laxman Sudhakar Padwale'
$ws.Cells.Item(9, 6).Formula = "'" + '033910330510'
$ws.Cells.Item(9, 8).Value = 'IPOS0000001'
$ws.Cells.Item(9, 9).Value = 'Corporate Office'
$ws.Cells.Item(9, 10).Value = 'Corporate Office'
$ws.Cells.Item(9, 11).Value = 'India'
$ws.Cells.Item(9, 12).Value = 'O|A|Laxman Sudhakar Padwale|033910330510||IPOS0000001|Corporate Office|Corporate Office|India'

# Row 10: Bhagirath Electricals & Borewell
$ws.Cells.Item(10, 1).Value = 4240035796
$ws.Cells.Item(10, 2).Value = 198305
$ws.Cells.Item(10, 3).Value = 'O'
$ws.Cells.Item(10, 4).Value = 'A'
$ws.Cells.Item(10, 5).Value = 'Bhagirath Electricals & Borewell'
$ws.Cells.Item(10, 6).Formula = "'" + '60056433432'
$ws.Cells.Item(10, 8).Value = 'MAHB0000429'
$ws.Cells.Item(10, 9).Value = 'Vikramgad'
$ws.Cells.Item(10, 10).Value = 'Vikramgad'
$ws.Cells.Item(10, 11).Value = 'India'
$ws.Cells.Item(10, 12).Value = 'O|A|Bhagirath Electricals & Borewell|60056433432||MAHB0000429|Vikramgad|Vikramgad|India'

# Row 11: Kore Mining & Crushing Pvt Ltd
$ws.Cells.Item(11, 1).Value = 4240035785
$ws.Cells.Item(11, 2).Value = 85374
$ws.Cells.Item(11, 3).Value = 'O'
$ws.Cells.Item(11, 4).Value = 'A'
$ws.Cells.Item(11, 5).Value = 'Kore Mining & Crushing Pvt Ltd'
$ws.Cells.Item(11, 6).Formula = "'" + '396100100001108'
$ws.Cells.Item(11, 8).Value = 'SRCB0000396'
$ws.Cells.Item(11, 9).Value = 'Boisar'
$ws.Cells.Item(11, 10).Value = 'Boisar'
$ws.Cells.Item(11, 11).Value = 'India'
$ws.Cells.Item(11, 12).Value = 'O|A|Kore Mining & Crushing Pvt Ltd|396100100001108||SRCB0000396|Boisar|Boisar|India'

# Row 12: Ganesh Ashok Sarode
$ws.Cells.Item(12, 1).Value = 4240035782
$ws.Cells.Item(12, 2).Value = 1533
$ws.Cells.Item(12, 3).Value = 'O'
$ws.Cells.Item(12, 4).Value = 'A'
$ws.Cells.Item(12, 5).Value = 'Ganesh Ashok Sarode'
$ws.Cells.Item(12, 6).Formula = "'" + '50100286551314'
$ws.Cells.Item(12, 8).Value = 'HDFC0002865'
$ws.Cells.Item(12, 9).Value = 'Jawhar'
$ws.Cells.Item(12, 10).Value = 'Jawhar'
$ws.Cells.Item(12, 11).Value = 'India'
$ws.Cells.Item(12, 12).Value = 'O|A|Ganesh Ashok Sarode|50100286551314||HDFC0002865|Jawhar|Jawhar|India'

# Row 13: Vijay Hardware Stores
$ws.Cells.Item(13, 1).Value = 4240035773
$ws.Cells.Item(13, 2).Value = 20945
$ws.Cells.Item(13, 3).Value = 'O'
$ws.Cells.Item(13, 4).Value = 'A'
$ws.Cells.Item(13, 5).Value = 'Vijay Hardware Stores'
$ws.Cells.Item(13, 6).Formula = "'" + '50200098814777'
$ws.Cells.Item(13, 8).Value = 'HDFC0007179'
$ws.Cells.Item(13, 9).Value = 'Manor'
$ws.Cells.Item(13, 10).Value = 'Manor'
$ws.Cells.Item(13, 11).Value = 'India'
$ws.Cells.Item(13, 12).Value = 'O|A|Vijay Hardware Stores|50200098814777||HDFC0007179|Manor|Manor|India'

# Row 14: Mohit Kumar
$ws.Cells.Item(14, 1).Value = 4240035770
$ws.Cells.Item(14, 2).Value = 15000
$ws.Cells.Item(14, 3).Value = 'S'
$ws.Cells.Item(14, 4).Value = 'A'
$ws.Cells.Item(14, 5).Value = 'Mohit Kumar'
$ws.Cells.Item(14, 6).Formula = "'" + '38796858574'
$ws.Cells.Item(14, 8).Value = 'SBIN0001143'
$ws.Cells.Item(14, 9).Value = 'SAUGOR UNIVERSITY'
$ws.Cells.Item(14, 10).Value = 'SAUGOR UNIVERSITY'
$ws.Cells.Item(14, 11).Value = 'India'
$ws.Cells.Item(14, 12).Value = 'S|A|Mohit Kumar|38796858574||SBIN0001143|SAUGOR UNIVERSITY|SAUGOR UNIVERSITY|India'

# Row 15: Ravi Shankar Lahange
$ws.Cells.Item(15, 1).Value = 4240035740
$ws.Cells.Item(15, 2).Value = 4036
$ws.Cells.Item(15, 3).Value = 'S'
$ws.Cells.Item(15, 4).Value = 'A'
$ws.Cells.Item(15, 5).Value = 'Ravi Shankar Lahange'
$ws.Cells.Item(15, 6).Formula = "'" + '35187115601'
$ws.Cells.Item(15, 8).Value = 'SBIN0007773'
$ws.Cells.Item(15, 9).Value = 'Gorhe'
$ws.Cells.Item(15, 10).Value = 'Gorhe'
$ws.Cells.Item(15, 11).Value = 'India'
$ws.Cells.Item(15, 12).Value = 'S|A|Ravi Shankar Lahange|35187115601||SBIN0007773|Gorhe|Gorhe|India'

# Row 16: Dipak Mahadu Gavit
$ws.Cells.Item(16, 1).Value = 4240035699
$ws.Cells.Item(16, 2).Value = 1783
$ws.Cells.Item(16, 3).Value = 'S'
$ws.Cells.Item(16, 4).Value = 'A'
$ws.Cells.Item(16, 5).Value = 'Dipak Mahadu Gavit'
$ws.Cells.Item(16, 6).Formula = "'" + '33523530337'
$ws.Cells.Item(16, 8).Value = 'SBIN0001050'
$ws.Cells.Item(16, 9).Value = 'JAWHAR'
$ws.Cells.Item(16, 10).Value = 'JAWHAR'
$ws.Cells.Item(16, 11).Value = 'India'
$ws.Cells.Item(16, 12).Value = 'S|A|Dipak Mahadu Gavit|33523530337||SBIN0001050|JAWHAR|JAWHAR|India'

# Row 17: Bharat Xerox And Stationary
$ws.Cells.Item(17, 1).Value = 4240035698
$ws.Cells.Item(17, 2).Value = 3084
$ws.Cells.Item(17, 3).Value = 'S'
$ws.Cells.Item(17, 4).Value = 'A'
$ws.Cells.Item(17, 5).Value = 'Bharat Xerox And Stationary'
$ws.Cells.Item(17, 6).Formula = "'" + '38166911649'
$ws.Cells.Item(17, 8).Value = 'SBIN0001050'
$ws.Cells.Item(17, 9).Value = 'Jawhar '
$ws.Cells.Item(17, 10).Value = 'Jawhar '
$ws.Cells.Item(17, 11).Value = 'India'
$ws.Cells.Item(17, 12).Value = 'S|A|Bharat Xerox And Stationary|38166911649||SBIN0001050|Jawhar |Jawhar |India'

# Row 18: Nilesh Vitthal Dhapshi
$ws.Cells.Item(18, 1).Value = 4240035427
$ws.Cells.Item(18, 2).Value = 3939
$ws.Cells.Item(18, 3).Value = 'S'
$ws.Cells.Item(18, 4).Value = 'A'
$ws.Cells.Item(18, 5).Value = 'Nilesh Vitthal Dhapshi'
$ws.Cells.Item(18, 6).Formula = "'" + '41568298765'
$ws.Cells.Item(18, 8).Value = 'SBIN0007773'
$ws.Cells.Item(18, 9).Value = 'Gorhe'
$ws.Cells.Item(18, 10).Value = 'Gorhe'
$ws.Cells.Item(18, 11).Value = 'India'
$ws.Cells.Item(18, 12).Value = 'S|A|Nilesh Vitthal Dhapshi|41568298765||SBIN0007773|Gorhe|Gorhe|India'

# Row 19: Madhavi Suryaji Desai
$ws.Cells.Item(19, 1).Value = 4240034980
$ws.Cells.Item(19, 2).Value = 17198
$ws.Cells.Item(19, 3).Value = 'O'
$ws.Cells.Item(19, 4).Value = 'A'
$ws.Cells.Item(19, 5).Value = 'Madhavi Suryaji Desai'
$ws.Cells.Item(19, 6).Formula = "'" + '60258069880'
$ws.Cells.Item(19, 8).Value = 'MAHB0000195'
$ws.Cells.Item(19, 9).Value = 'Palghar'
$ws.Cells.Item(19, 10).Value = 'Palghar'
$ws.Cells.Item(19, 11).Value = 'India'
$ws.Cells.Item(19, 12).Value = 'O|A|Madhavi Suryaji Desai|60258069880||MAHB0000195|Palghar|Palghar|India'

# Row 20: Ankush Yadav
$ws.Cells.Item(20, 1).Value = 4240036101
$ws.Cells.Item(20, 2).Value = 3589
$ws.Cells.Item(20, 3).Value = 'O'
$ws.Cells.Item(20, 4).Value = 'A'
$ws.Cells.Item(20, 5).Value = 'Ankush Yadav'
$ws.Cells.Item(20, 6).Formula = "'" + '681518110001008'
$ws.Cells.Item(20, 8).Value = 'BKID0006815'
$ws.Cells.Item(20, 9).Value = 'Subhanikheda'
$ws.Cells.Item(20, 10).Value = 'Subhanikheda'
$ws.Cells.Item(20, 11).Value = 'India'
$ws.Cells.Item(20, 12).Value = 'O|A|Ankush Yadav|681518110001008||BKID0006815|Subhanikheda|Subhanikheda|India'

# Row 21: TOTAL
$ws.Cells.Item(21, 1).Value = "TOTAL"
$ws.Cells.Item(21, 2).Value = 404985
